$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.246.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.430.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.24%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "489.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.611"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +18.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.434.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.100"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.333"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.69%  "

$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.842.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.265.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.46%  "

$ws.Range("E17").Value = "  -3.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.430.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "57.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.11%  "

$ws.Range("E25").Value = "  -1.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.32%  "

$ws.Range("E27").Value = "  -3.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.524.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0786"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.25%  "

$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.17%  "

$ws.Range("E34").Value = "  -1.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  -0.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.82%  "

$ws.Range("E38").Value = "  -10.20%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "285.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.12%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.101"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.73%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "33.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.33%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.98%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.596"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0532"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0228"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.911.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.69%  "
